$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.921.18'
$ws.Range("E2").Value = '  +4.30%  '

# Row 3
$ws.Range("D3").Value = '2.275.34'
$ws.Range("E3").Value = '  +4.87%  '

# Row 4
$ws.Range("E4").Value = '  +0.41%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.36'
$ws.Range("E5").Value = '  +1.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.636'
$ws.Range("E6").Value = '  +3.92%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.89'
$ws.Range("E7").Value = '  +9.58%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.664'
$ws.Range("E9").Value = '  +18.73%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.89'
$ws.Range("E10").Value = '  +10.14%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.62'
$ws.Range("E11").Value = '  +0.07%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0970'
$ws.Range("E12").Value = '  +4.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.45'
$ws.Range("E13").Value = '  +9.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.104'
$ws.Range("E14").Value = '  +0.86%  '

# Row 15
$ws.Range("D15").Value = '2.616.76'
$ws.Range("E15").Value = '  +5.07%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.89'
$ws.Range("E16").Value = '  +4.57%  '

# Row 17
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.888'
$ws.Range("E17").Value = '  +5.02%  '

# Row 18
$ws.Range("D18").Value = '2.283.94'
$ws.Range("E18").Value = '  +5.71%  '

# Row 19
$ws.Range("D19").Value = '42.862.24'
$ws.Range("E19").Value = '  +4.46%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000100'
$ws.Range("E20").Value = '  +7.66%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("E21").Value = '  +4.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.20'
$ws.Range("E22").Value = '  +2.69%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.10'
$ws.Range("E23").Value = '  +3.29%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  +3.32%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.12'
$ws.Range("E25").Value = '  +7.85%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.08%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("E27").Value = '  +2.19%  '

# Row 28
$ws.Range("E28").Value = '  +1.14%  '

# Row 29
$ws.Range("E29").Value = '  -1.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +4.73%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.54'
$ws.Range("E31").Value = '  +0.06%  '

# Row 32
$ws.Range("E32").Value = '  +4.60%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.58'
$ws.Range("E33").Value = '  +17.16%  '

# Row 34
$ws.Range("E34").Value = '  +5.91%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0805'
$ws.Range("E35").Value = '  +8.82%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.57'
$ws.Range("E36").Value = '  +29.84%  '

# Row 37
$ws.Range("E37").Value = '  +4.32%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.46'
$ws.Range("E38").Value = '  +12.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.77'
$ws.Range("E39").Value = '  +5.25%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0316'
$ws.Range("E40").Value = '  +5.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +7.33%  '

# Row 42
$ws.Range("E42").Value = '  +12.67%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.82'
$ws.Range("E43").Value = '  +6.92%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.208'
$ws.Range("E44").Value = '  +9.12%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.28'
$ws.Range("E45").Value = '  +10.09%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.23'
$ws.Range("E46").Value = '  +3.20%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.89'
$ws.Range("E47").Value = '  +1.24%  '

# Row 48
$ws.Range("E48").Value = '  +3.99%  '

# Row 49
$ws.Range("E49").Value = '  +0.12%  '

# Row 50
$ws.Range("E50").Value = '  +3.07%  '

# Row 51
$ws.Range("E51").Value = '  +4.73%  '
